# Feria Lagunitas de Puerto Montt - Zanahoria: add a new weekly record.
# A new row is inserted above row 169, shifting the existing rows 169:241
# down to 170:242, and the new row 169 is populated with the latest
# price observation (date 2021-11-10, serial 44510).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 169 (pushes rows 169-241 down to 170-242).
$ws.Rows(169).Insert()

# Fill in the new row 169 with the new record's data.
$ws.Cells.Item(169, 1).Value  = 4
$ws.Cells.Item(169, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(169, 3).Value  = "Los Lagos"
$ws.Cells.Item(169, 4).Value  = 44510
$ws.Cells.Item(169, 5).Value  = 10
$ws.Cells.Item(169, 6).Value  = 100114013
$ws.Cells.Item(169, 7).Value  = "Zanahoria"
$ws.Cells.Item(169, 8).Value  = "Sin especificar"
$ws.Cells.Item(169, 9).Value  = "Primera"
$ws.Cells.Item(169, 10).Value = 150
$ws.Cells.Item(169, 11).Value = 12000
$ws.Cells.Item(169, 12).Value = 12000
$ws.Cells.Item(169, 13).Value = 12000
$ws.Cells.Item(169, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(169, 15).Value = "Región de Ñuble"
$ws.Cells.Item(169, 16).Value = 600
$ws.Cells.Item(169, 17).Value = 20
$ws.Cells.Item(169, 18).Value = "Hortaliza"
